$wb = $excel.ActiveWorkbook

# Locate the "Italy" worksheet, which acts as the template for the new "Spain" sheet.
$italy = $wb.Worksheets.Item("Italy")

# Copy the Italy sheet to create the new Spain sheet, placing it right after Italy.
$italy.Copy([System.Reflection.Missing]::Value, $italy)
$spain = $wb.Worksheets.Item($italy.Index + 1)
$spain.Name = "Spain"

# Update the market-specific values on the new sheet.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2037 "

# Restore Italy's selection to the full used range (no longer the active tab).
$italy.Activate()
$italy.Range("A1:D15").Select()

# Make Spain the active sheet/tab and set its selection (B4:B5, active cell B5).
$spain.Activate()
$spain.Range("B4:B5").Select()
